$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The model-metrics table is fully rewritten: 25 models (_0.._24) instead of
# 20 (_0.._19), re-ranked by this training run (rstd improved for 4 more
# models - "4 mdelo melhores rstds"). Clear the numeric columns (B:I) across
# the old + new extent; column A is left alone for the rows that already
# existed so their existing label style (bold/bordered/centered, style index
# 1) is preserved untouched.
$ws.Range("B2:I30").Clear()

# Header row (B1:I1) is unchanged, but (re)write defensively.
$headers = @("r2","r2_test","r2_val","r2_vt","mse","mse_test","mse_val","mse_vt")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# New model metrics data (rows 2-26): label, r2, r2_test, r2_val, r2_vt, mse,
# mse_test, mse_val, mse_vt.
$data = @(
    ,@("model_7_9_24", 0.8051527574306611, 0.3571871387255203, -3.414186181067846, -0.0950500691864502, 0.2156384885311127, 0.6085723638534546, 0.7437629699707031, 0.6721917986869812)
    ,@("model_7_9_23", 0.8059775252917172, 0.3611886883211213, -3.381981137552415, -0.08762251096688445, 0.2147257030010223, 0.6047839522361755, 0.7383365631103516, 0.6676324605941772)
    ,@("model_7_9_22", 0.8067786356822158, 0.3652607565990997, -3.350864644433456, -0.0802782813183851, 0.213839128613472, 0.6009288430213928, 0.7330936193466187, 0.6631242036819458)
    ,@("model_7_9_21", 0.8075879361699952, 0.369249891183979, -3.31746852490921, -0.07270708278244387, 0.2129434496164322, 0.5971521735191345, 0.7274665832519531, 0.6584766507148743)
    ,@("model_7_9_20", 0.8084681234697515, 0.3737561330882341, -3.28339345358619, -0.06462603265679978, 0.2119693458080292, 0.5928859710693359, 0.7217251062393188, 0.6535161137580872)
    ,@("model_7_9_19", 0.8093878739994109, 0.3781107403989143, -3.245661194444701, -0.05619699178946203, 0.2109514474868774, 0.5887632966041565, 0.7153674960136414, 0.6483420133590698)
    ,@("model_7_9_18", 0.8102836627551147, 0.3828586816190225, -3.210487823013826, -0.04777721437983673, 0.209960088133812, 0.5842682719230652, 0.7094410061836243, 0.6431735157966614)
    ,@("model_7_9_17", 0.8111397944438427, 0.3875931499582037, -3.175888168681357, -0.03944199835004714, 0.2090125977993011, 0.5797860026359558, 0.7036111354827881, 0.6380570530891418)
    ,@("model_7_9_16", 0.8120333461211127, 0.3922411821162277, -3.138048322955895, -0.03075728167006431, 0.2080236822366714, 0.575385570526123, 0.6972353458404541, 0.6327259540557861)
    ,@("model_7_9_15", 0.8129445311470461, 0.3969942450583754, -3.09937117833863, -0.02188083745649272, 0.2070153057575226, 0.5708855986595154, 0.6907185316085815, 0.627277135848999)
    ,@("model_7_9_14", 0.8137078040993939, 0.4009877974847715, -3.060737272901118, -0.01362988729968406, 0.2061705589294434, 0.5671048760414124, 0.6842089295387268, 0.6222124099731445)
    ,@("model_7_9_13", 0.814544113019412, 0.4060163579901461, -3.024955661981612, -0.004901226901900868, 0.2052450180053711, 0.5623441338539124, 0.6781800389289856, 0.6168543100357056)
    ,@("model_7_9_12", 0.8153089468311172, 0.4106813086730439, -2.988251207413242, 0.003646604703333001, 0.2043985575437546, 0.5579277276992798, 0.6719955205917358, 0.6116072535514832)
    ,@("model_7_9_11", 0.8160603455728173, 0.4161427624766562, -2.956438451727058, 0.01221548384536386, 0.2035669982433319, 0.5527570843696594, 0.6666352152824402, 0.6063473224639893)
    ,@("model_7_9_10", 0.8167924148382092, 0.4211846534008626, -2.921935748942603, 0.02079054209161513, 0.2027567923069, 0.5479838252067566, 0.660821795463562, 0.6010834574699402)
    ,@("model_7_9_9", 0.8175460424441786, 0.4264303189421452, -2.886235028068094, 0.02968380691646466, 0.2019227594137192, 0.543017566204071, 0.6548064947128296, 0.5956244468688965)
    ,@("model_7_9_8", 0.8181541341143521, 0.4319044790902338, -2.856603188249268, 0.03798165264255637, 0.2012497633695602, 0.5378349423408508, 0.6498136520385742, 0.5905308723449707)
    ,@("model_7_9_7", 0.8187747907233892, 0.4376679490047077, -2.827230820904878, 0.04648225729135669, 0.200562909245491, 0.5323784947395325, 0.6448646187782288, 0.5853127837181091)
    ,@("model_7_9_6", 0.819294623741498, 0.4426451514574946, -2.795692205420966, 0.05462016926893531, 0.1999876350164413, 0.527666449546814, 0.639550507068634, 0.5803173780441284)
    ,@("model_7_9_5", 0.8197217330019009, 0.4477244674822253, -2.767077397166287, 0.06246258941991101, 0.1995149105787277, 0.5228577256202698, 0.6347291469573975, 0.5755033493041992)
    ,@("model_7_9_4", 0.8199889714809315, 0.4527281690723072, -2.744374830428599, 0.06948024541546527, 0.1992191523313522, 0.5181205272674561, 0.6309038400650024, 0.5711955428123474)
    ,@("model_7_9_3", 0.8201706992472464, 0.4579003664384543, -2.724744402352817, 0.07623974519756527, 0.1990180313587189, 0.5132238268852234, 0.6275962591171265, 0.5670462846755981)
    ,@("model_7_9_2", 0.8202694615576511, 0.4627862888377087, -2.705261233096862, 0.0827463814110232, 0.198908731341362, 0.5085981488227844, 0.6243134737014771, 0.563052237033844)
    ,@("model_7_9_0", 0.8203433966964468, 0.4719407136084411, -2.660441333820163, 0.09601129379997952, 0.1988269090652466, 0.4999313950538635, 0.6167616248130798, 0.5549095273017883)
    ,@("model_7_9_1", 0.8204704538135593, 0.4676024849160307, -2.678108489783365, 0.09018606877622792, 0.1986862868070602, 0.5040385723114014, 0.6197383999824524, 0.5584853887557983)
)

$startRow = 2
$lastRow = $startRow + $data.Length - 1
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $startRow + $r
    $ws.Cells.Item($excelRow, 1).Value = $rowVals[0]
    for ($c = 1; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowVals[$c]
    }
}

# Extend the label style (bold, thin border, centered/top) used on column A
# down to the newly added rows (22-26), by copying the format from an
# already-styled label cell. This reuses the existing style record (index 1)
# instead of minting a new, functionally-identical one.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A22:A" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

